# Update "想去人数" (interest count) values in column F across sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1191
$ws1.Range("F4").Value = 48
$ws1.Range("F5").Value = 1312
$ws1.Range("F6").Value = 1696
$ws1.Range("F7").Value = 6214
$ws1.Range("F9").Value = 1814
$ws1.Range("F10").Value = 481
$ws1.Range("F12").Value = 15
$ws1.Range("F15").Value = 21
$ws1.Range("F16").Value = 6908
$ws1.Range("F18").Value = 53
$ws1.Range("F19").Value = 164
$ws1.Range("F21").Value = 1703
$ws1.Range("F23").Value = 14
$ws1.Range("F24").Value = 43
$ws1.Range("F25").Value = 161
$ws1.Range("F26").Value = 1576
$ws1.Range("F27").Value = 756
$ws1.Range("F28").Value = 313
$ws1.Range("F31").Value = 52
$ws1.Range("F33").Value = 3891

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 345
$ws2.Range("F5").Value = 203
$ws2.Range("F8").Value = 444
$ws2.Range("F19").Value = 5

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 9527
$ws3.Range("F3").Value = 2260
$ws3.Range("F5").Value = 242

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9527
$ws4.Range("F3").Value = 2260
$ws4.Range("F5").Value = 1191
$ws4.Range("F7").Value = 48
$ws4.Range("F9").Value = 345
$ws4.Range("F10").Value = 1312
$ws4.Range("F11").Value = 242
$ws4.Range("F12").Value = 1696
$ws4.Range("F13").Value = 6214
$ws4.Range("F15").Value = 1814
$ws4.Range("F18").Value = 481
$ws4.Range("F20").Value = 15
$ws4.Range("F23").Value = 6908
$ws4.Range("F25").Value = 53
$ws4.Range("F26").Value = 164
$ws4.Range("F28").Value = 1703
$ws4.Range("F30").Value = 14
$ws4.Range("F31").Value = 43
$ws4.Range("F32").Value = 161
$ws4.Range("F33").Value = 1576
$ws4.Range("F34").Value = 756
$ws4.Range("F36").Value = 313
$ws4.Range("F44").Value = 5
$ws4.Range("F45").Value = 3891
